$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve the footer ("data source") row style before it gets overwritten ---
$ws.Range("A3126").Copy()
$ws.Range("A3165").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Copy the numeric/date cell formatting down across the rows that will now hold data ---
$ws.Range("A3122:D3122").Copy()
$ws.Range("A3124:D3162").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Match the row height used by the other data rows ---
$ws.Range("A3124:A3162").RowHeight = 14

# --- Fill in the previously-missing spot price for 2024-01-19 ---
$ws.Cells.Item(3123, 2).Value = 2212.9

# --- Append the new daily observations (coke spot price, futures close, futures settle) ---
$ws.Cells.Item(3124, 1).Value = 45313
$ws.Cells.Item(3124, 2).Value = 2212.9
$ws.Cells.Item(3124, 3).Value = 2454
$ws.Cells.Item(3124, 4).Value = 2464
$ws.Cells.Item(3125, 1).Value = 45314
$ws.Cells.Item(3125, 2).Value = 2212.9
$ws.Cells.Item(3125, 3).Value = 2482
$ws.Cells.Item(3125, 4).Value = 2472.5
$ws.Cells.Item(3126, 1).Value = 45315
$ws.Cells.Item(3126, 2).Value = 2212.9
$ws.Cells.Item(3126, 3).Value = 2504
$ws.Cells.Item(3126, 4).Value = 2486
$ws.Cells.Item(3127, 1).Value = 45316
$ws.Cells.Item(3127, 2).Value = 2212.9
$ws.Cells.Item(3127, 3).Value = 2493.5
$ws.Cells.Item(3127, 4).Value = 2510
$ws.Cells.Item(3128, 1).Value = 45317
$ws.Cells.Item(3128, 2).Value = 2212.9
$ws.Cells.Item(3128, 3).Value = 2459.5
$ws.Cells.Item(3128, 4).Value = 2479
$ws.Cells.Item(3129, 1).Value = 45320
$ws.Cells.Item(3129, 2).Value = 2212.9
$ws.Cells.Item(3129, 3).Value = 2444
$ws.Cells.Item(3129, 4).Value = 2460.5
$ws.Cells.Item(3130, 1).Value = 45321
$ws.Cells.Item(3130, 2).Value = 2212.9
$ws.Cells.Item(3130, 3).Value = 2398
$ws.Cells.Item(3130, 4).Value = 2433
$ws.Cells.Item(3131, 1).Value = 45322
$ws.Cells.Item(3131, 2).Value = 2212.9
$ws.Cells.Item(3131, 3).Value = 2346
$ws.Cells.Item(3131, 4).Value = 2373
$ws.Cells.Item(3132, 1).Value = 45323
$ws.Cells.Item(3132, 2).Value = 2212.9
$ws.Cells.Item(3132, 3).Value = 2348.5
$ws.Cells.Item(3132, 4).Value = 2347
$ws.Cells.Item(3133, 1).Value = 45324
$ws.Cells.Item(3133, 2).Value = 2212.9
$ws.Cells.Item(3133, 3).Value = 2328
$ws.Cells.Item(3133, 4).Value = 2337
$ws.Cells.Item(3134, 1).Value = 45327
$ws.Cells.Item(3134, 2).Value = 2212.9
$ws.Cells.Item(3134, 3).Value = 2332.5
$ws.Cells.Item(3134, 4).Value = 2337.5
$ws.Cells.Item(3135, 1).Value = 45328
$ws.Cells.Item(3135, 2).Value = 2212.9
$ws.Cells.Item(3135, 3).Value = 2330
$ws.Cells.Item(3135, 4).Value = 2316
$ws.Cells.Item(3136, 1).Value = 45329
$ws.Cells.Item(3136, 2).Value = 2212.9
$ws.Cells.Item(3136, 3).Value = 2317.5
$ws.Cells.Item(3136, 4).Value = 2320
$ws.Cells.Item(3137, 1).Value = 45330
$ws.Cells.Item(3137, 2).Value = 2212.9
$ws.Cells.Item(3137, 3).Value = 2378
$ws.Cells.Item(3137, 4).Value = 2346
$ws.Cells.Item(3138, 1).Value = 45341
$ws.Cells.Item(3138, 2).Value = 2112.9
$ws.Cells.Item(3138, 3).Value = 2271
$ws.Cells.Item(3138, 4).Value = 2307.5
$ws.Cells.Item(3139, 1).Value = 45342
$ws.Cells.Item(3139, 2).Value = 2112.9
$ws.Cells.Item(3139, 3).Value = 2255.5
$ws.Cells.Item(3139, 4).Value = 2269.5
$ws.Cells.Item(3140, 1).Value = 45343
$ws.Cells.Item(3140, 2).Value = 2112.9
$ws.Cells.Item(3140, 3).Value = 2355.5
$ws.Cells.Item(3140, 4).Value = 2334
$ws.Cells.Item(3141, 1).Value = 45344
$ws.Cells.Item(3141, 2).Value = 2112.9
$ws.Cells.Item(3141, 3).Value = 2403.5
$ws.Cells.Item(3141, 4).Value = 2381
$ws.Cells.Item(3142, 1).Value = 45345
$ws.Cells.Item(3142, 2).Value = 2112.9
$ws.Cells.Item(3142, 3).Value = 2370.5
$ws.Cells.Item(3142, 4).Value = 2383.5
$ws.Cells.Item(3143, 1).Value = 45348
$ws.Cells.Item(3143, 2).Value = 2112.9
$ws.Cells.Item(3143, 3).Value = 2314
$ws.Cells.Item(3143, 4).Value = 2344.5
$ws.Cells.Item(3144, 1).Value = 45349
$ws.Cells.Item(3144, 2).Value = 2012.9
$ws.Cells.Item(3144, 3).Value = 2394
$ws.Cells.Item(3144, 4).Value = 2357.5
$ws.Cells.Item(3145, 1).Value = 45350
$ws.Cells.Item(3145, 2).Value = 2012.9
$ws.Cells.Item(3145, 3).Value = 2389
$ws.Cells.Item(3145, 4).Value = 2382
$ws.Cells.Item(3146, 1).Value = 45351
$ws.Cells.Item(3146, 2).Value = 2012.9
$ws.Cells.Item(3146, 3).Value = 2381.5
$ws.Cells.Item(3146, 4).Value = 2376.5
$ws.Cells.Item(3147, 1).Value = 45352
$ws.Cells.Item(3147, 2).Value = 2012.9
$ws.Cells.Item(3147, 3).Value = 2361.5
$ws.Cells.Item(3147, 4).Value = 2394.5
$ws.Cells.Item(3148, 1).Value = 45355
$ws.Cells.Item(3148, 2).Value = 2012.9
$ws.Cells.Item(3148, 3).Value = 2361
$ws.Cells.Item(3148, 4).Value = 2347
$ws.Cells.Item(3149, 1).Value = 45356
$ws.Cells.Item(3149, 2).Value = 2012.9
$ws.Cells.Item(3149, 3).Value = 2326.5
$ws.Cells.Item(3149, 4).Value = 2337.5
$ws.Cells.Item(3150, 1).Value = 45357
$ws.Cells.Item(3150, 2).Value = 2012.9
$ws.Cells.Item(3150, 3).Value = 2291
$ws.Cells.Item(3150, 4).Value = 2304.5
$ws.Cells.Item(3151, 1).Value = 45358
$ws.Cells.Item(3151, 2).Value = 2012.9
$ws.Cells.Item(3151, 3).Value = 2300
$ws.Cells.Item(3151, 4).Value = 2304
$ws.Cells.Item(3152, 1).Value = 45359
$ws.Cells.Item(3152, 2).Value = 2012.9
$ws.Cells.Item(3152, 3).Value = 2292
$ws.Cells.Item(3152, 4).Value = 2283
$ws.Cells.Item(3153, 1).Value = 45362
$ws.Cells.Item(3153, 2).Value = 2012.9
$ws.Cells.Item(3153, 3).Value = 2236.5
$ws.Cells.Item(3153, 4).Value = 2251
$ws.Cells.Item(3154, 1).Value = 45363
$ws.Cells.Item(3154, 2).Value = 1912.9
$ws.Cells.Item(3154, 3).Value = 2244
$ws.Cells.Item(3154, 4).Value = 2248
$ws.Cells.Item(3155, 1).Value = 45364
$ws.Cells.Item(3155, 2).Value = 1912.9
$ws.Cells.Item(3155, 3).Value = 2194.5
$ws.Cells.Item(3155, 4).Value = 2207.5
$ws.Cells.Item(3156, 1).Value = 45365
$ws.Cells.Item(3156, 2).Value = 1912.9
$ws.Cells.Item(3156, 3).Value = 2168.5
$ws.Cells.Item(3156, 4).Value = 2178.5
$ws.Cells.Item(3157, 1).Value = 45366
$ws.Cells.Item(3157, 2).Value = 1912.9
$ws.Cells.Item(3157, 3).Value = 2141
$ws.Cells.Item(3157, 4).Value = 2147
$ws.Cells.Item(3158, 1).Value = 45369
$ws.Cells.Item(3158, 2).Value = 1912.9
$ws.Cells.Item(3158, 3).Value = 2147
$ws.Cells.Item(3158, 4).Value = 2126.5
$ws.Cells.Item(3159, 1).Value = 45370
$ws.Cells.Item(3159, 2).Value = 1912.9
$ws.Cells.Item(3159, 3).Value = 2179.5
$ws.Cells.Item(3159, 4).Value = 2159.5
$ws.Cells.Item(3160, 1).Value = 45371
$ws.Cells.Item(3160, 2).Value = 1812.9
$ws.Cells.Item(3160, 3).Value = 2175
$ws.Cells.Item(3160, 4).Value = 2185
$ws.Cells.Item(3161, 1).Value = 45372
$ws.Cells.Item(3161, 2).Value = 1812.9
$ws.Cells.Item(3161, 3).Value = 2186
$ws.Cells.Item(3161, 4).Value = 2172.5
$ws.Cells.Item(3162, 1).Value = 45373
$ws.Cells.Item(3162, 3).Value = 2185.5
$ws.Cells.Item(3162, 4).Value = 2178

# --- Move the "data source" footer note down below the newly appended rows ---
$ws.Range("A3165").Value = "数据来源：东方财富Choice数据"

# --- Reflect the new active cell / selection state ---
$ws.Range("D3165").Select()
